$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "; Field Type should be..." -> "; field Type should be..." (lower-case
#    the "F") inside the "Create a new calculated field..." bullet.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("; Field Type should be", $true, $false, $false, $false, $false, $true, 1, $false, "; field Type should be", 2)

# ---------------------------------------------------------------------------
# 2. Insert a brand new bullet right after the "spatially constrained
#    multivariate clustering tool" bullet, describing the join of the
#    [#]_Clusters table with the VetData_Centroids table on XCord.
# ---------------------------------------------------------------------------
$pClusterTool = $d.Paragraphs.Item(7)
$pClusterTool.Range.InsertParagraphAfter()

$pNewJoin = $d.Paragraphs.Item(8)
$pNewJoin.Range.Text = "Join the [State]_Census_Blocks_VetData_[#]_Clusters table with the [State]_Census_Blocks_VetData_Centroids, using field XCord and using the former table as the primary."

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the end of the "Use the dissolve
#    tool...VetSums" bullet (paragraph 5) to the end of the bullet we just
#    inserted (paragraph 8).
# ---------------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$pTarget = $d.Paragraphs.Item(8)
$rngEnd = $pTarget.Range.Duplicate
$rngEnd.MoveEnd(1, -1)
$rngEnd.Collapse(0)
$bmPos = $rngEnd.Start

# Work around a collapsed-range-at-paragraph-boundary quirk: temporarily
# insert a placeholder character, anchor the bookmark to it, then delete the
# placeholder so the bookmark collapses back to a true zero-width mark.
$rngEnd.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $rngEnd)
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
